$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-1h (E) columns for rows with new market data
$ws.Range("D2").Value = "30.646.77"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "2.119.93"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.96%  "
$ws.Range("D5").Value = "337.84"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("D7").Value = "0.5257"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("D8").Value = "0.4558"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "54.47"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").Value = "0.09125"
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("D11").Value = "1.177"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "24.51"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").Value = "2.124.72"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "6.868"
$ws.Range("D15").Value = "8.146"
$ws.Range("E15").Value = "  +5.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001175"
$ws.Range("E16").Value = "  +4.66%  "
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D20").Value = "19.48"
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "6.322"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").Value = "30.722.17"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "12.91"
$ws.Range("E24").Value = "  +4.90%  "
$ws.Range("D25").Value = "2.364"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D26").Value = "2.372.49"
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("D27").Value = "22.41"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").Value = "164.88"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").Value = "2.558"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Value = "134.85"
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.210"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("D32").Value = "0.1075"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "1.656"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.370"
$ws.Range("E34").Value = "  +3.31%  "
$ws.Range("D35").Value = "3.945"
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("D36").Value = "10.71"
$ws.Range("E36").Value = "  +5.56%  "
$ws.Range("D37").Value = "5.874"
$ws.Range("E37").Value = "  +7.43%  "
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06870"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("E40").Value = "  +3.19%  "
$ws.Range("D41").Value = "12.67"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").Value = "0.6922"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").Value = "1.259"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").Value = "14.99"
$ws.Range("E44").Value = "  +6.89%  "
$ws.Range("D45").Value = "0.6507"
$ws.Range("E45").Value = "  +2.72%  "
$ws.Range("D46").Value = "2.318"
$ws.Range("E46").Value = "  +5.16%  "
$ws.Range("E47").Value = "  +23.43%  "
$ws.Range("D48").Value = "3.698"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("D49").Value = "1.257"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").Value = "83.53"
$ws.Range("E50").Value = "  +2.18%  "

# Row 51: WEMIXTOKEN replaced by Cronos
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07292"
$ws.Range("E51").Value = "  +3.50%  "
